# The commit removes the yellow highlight from the "Showing a special
# character (stick figure) on the LED screen" bullet, and relocates the
# "_GoBack" bookmark from the blank paragraph right after
# "**have input stored for player one DONE" down to the very last
# (blank, centered) paragraph of the document.

$d = $word.ActiveDocument

# 1) Strip the yellow highlight from the "Showing a special character..."
#    bullet. Setting Font.HighlightColorIndex on the paragraph's Range
#    clears the highlight both on the run and on the paragraph mark.
$highlightPara = $d.Paragraphs.Item(29)
$highlightPara.Range.Font.HighlightColorIndex = 0

# 2) Move the "_GoBack" bookmark to the last paragraph in the document.
#    Adding a bookmark with a name that already exists relocates it, so
#    this both removes it from its old spot and creates it at the new one.
$lastPara = $d.Paragraphs.Item(30)
$d.Bookmarks.Add("_GoBack", $lastPara.Range)
